# Adds a second sentence/run (the "Lorem Ipsum" filler paragraph) to the
# existing "Le scrivo questa brevissima mail..." paragraph, giving both the
# original sentence and the new sentence explicit run/paragraph-mark
# character formatting (Open Sans / 10.5pt / black), and nudges a couple of
# related paragraph-level settings on the "Normal" style to match.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Locate the target paragraph robustly (don't assume a fixed index).
# ---------------------------------------------------------------------
$search = $d.Content
$found = $search.Find.Execute(
    "Le scrivo questa brevissima mail di test per il progetto",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate target paragraph"
}

$para = $search.Paragraphs(1)
$pRange = $para.Range

$lorem = "Contrary to popular belief, Lorem Ipsum is not simply random text. It has roots in a piece of classical Latin literature from 45 BC, making it over 2000 years old. Richard McClintock, a Latin professor at Hampden-Sydney College in Virginia, looked up one of the more obscure Latin words, consectetur, from a Lorem Ipsum passage, and going through the cites of the word in classical literature, discovered the undoubtable source. Lorem Ipsum comes from sections 1.10.32 and 1.10.33 of `"de Finibus Bonorum et Malorum`" (The Extremes of Good and Evil) by Cicero, written in 45 BC. This book is a treatise on the theory of ethics, very popular during the Renaissance. The first line of Lorem Ipsum, `"Lorem ipsum dolor sit amet..`", comes from a line in section 1.10.32."

$runProps = '<w:rFonts w:ascii="Open Sans;Arial;sans-serif" w:hAnsi="Open Sans;Arial;sans-serif"/><w:b w:val="false"/><w:i w:val="false"/><w:caps w:val="false"/><w:smallCaps w:val="false"/><w:color w:val="000000"/><w:spacing w:val="0"/><w:sz w:val="21"/>'

$xml = '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:bidi w:val="0"/><w:jc w:val="left"/><w:rPr>' + $runProps + '</w:rPr></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">Le scrivo questa brevissima mail di test per il progetto ${nome_progetto}. </w:t></w:r><w:r><w:rPr>' + $runProps + '</w:rPr><w:t>' + $lorem + '</w:t></w:r></w:p>'

$pRange.InsertXML($xml)

# ---------------------------------------------------------------------
# 2. Paragraph-format tweaks on the "Normal" style (spacing/justification
#    and hyphenation) that accompanied the content change.
# ---------------------------------------------------------------------
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.SpaceBefore = 0
$normal.ParagraphFormat.SpaceAfter = 0
$normal.ParagraphFormat.Alignment = 0
$normal.ParagraphFormat.Hyphenation = $false

Write-Host "Edit applied"
